# Applies refreshed market-board pricing values to the Leve profit calculation
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), as produced by the scheduled
# data-refresh runner. Only the affected numeric cells are updated in place.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 177833.5
$ws.Cells.Item(70, 9).Value = 339332.66
$ws.Cells.Item(70, 11).Value = 1017997.98
$ws.Cells.Item(70, 13).Value = -1017727.98
$ws.Cells.Item(73, 8).Value = 177833.5
$ws.Cells.Item(73, 9).Value = 339332.66
$ws.Cells.Item(73, 11).Value = 1017997.98
$ws.Cells.Item(73, 13).Value = -1017061.98
$ws.Cells.Item(95, 8).Value = 65674.664
$ws.Cells.Item(95, 10).Value = 65674.664
$ws.Cells.Item(95, 12).Value = 65674.664
$ws.Cells.Item(95, 14).Value = -71166.664
$ws.Cells.Item(112, 8).Value = 1265.6207
$ws.Cells.Item(112, 10).Value = 1359.75
$ws.Cells.Item(112, 12).Value = 4079.25
$ws.Cells.Item(112, 14).Value = -6295.25
$ws.Cells.Item(141, 8).Value = 5116.6665
$ws.Cells.Item(141, 9).Value = 4350
$ws.Cells.Item(141, 11).Value = 13050
$ws.Cells.Item(141, 13).Value = -7870

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 222.33333
$ws.Cells.Item(5, 9).Value = 110.166664
$ws.Cells.Item(5, 10).Value = 334.5
$ws.Cells.Item(5, 11).Value = 110.166664
$ws.Cells.Item(5, 12).Value = 334.5
$ws.Cells.Item(5, 13).Value = 1.833336000000003
$ws.Cells.Item(5, 14).Value = -558.5
$ws.Cells.Item(32, 8).Value = 2690.3972
$ws.Cells.Item(32, 9).Value = 2049.261
$ws.Cells.Item(32, 11).Value = 2049.261
$ws.Cells.Item(32, 13).Value = -1762.261
$ws.Cells.Item(61, 8).Value = 4540.4287
$ws.Cells.Item(61, 9).Value = 3561.8262
$ws.Cells.Item(61, 11).Value = 3561.8262
$ws.Cells.Item(61, 13).Value = -3349.8262
$ws.Cells.Item(74, 8).Value = 13892209
$ws.Cells.Item(74, 9).Value = 18520438
$ws.Cells.Item(74, 11).Value = 18520438
$ws.Cells.Item(74, 13).Value = -18519564
$ws.Cells.Item(77, 8).Value = 13892209
$ws.Cells.Item(77, 9).Value = 18520438
$ws.Cells.Item(77, 11).Value = 92602190
$ws.Cells.Item(77, 13).Value = -92597822
$ws.Cells.Item(136, 8).Value = 4540.4287
$ws.Cells.Item(136, 9).Value = 3561.8262
$ws.Cells.Item(136, 11).Value = 10685.4786
$ws.Cells.Item(136, 13).Value = -8135.4786

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 222.33333
$ws.Cells.Item(4, 9).Value = 110.166664
$ws.Cells.Item(4, 10).Value = 334.5
$ws.Cells.Item(4, 11).Value = 110.166664
$ws.Cells.Item(4, 12).Value = 334.5
$ws.Cells.Item(4, 13).Value = 4.833336000000003
$ws.Cells.Item(4, 14).Value = -564.5
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).Value = ""
$ws.Cells.Item(95, 8).Value = 50001
$ws.Cells.Item(95, 10).Value = 50001
$ws.Cells.Item(95, 12).Value = 50001
$ws.Cells.Item(95, 14).Value = -55493
$ws.Cells.Item(105, 8).Value = 41213
$ws.Cells.Item(105, 9).Value = 42398.4
$ws.Cells.Item(105, 11).Value = 42398.4
$ws.Cells.Item(105, 13).Value = -40651.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 241.94737
$ws.Cells.Item(7, 9).Value = 32.5
$ws.Cells.Item(7, 11).Value = 32.5
$ws.Cells.Item(7, 13).Value = 80.5
$ws.Cells.Item(22, 8).Value = 2570.2856
$ws.Cells.Item(22, 9).Value = 499.25
$ws.Cells.Item(22, 10).Value = 5331.6665
$ws.Cells.Item(22, 11).Value = 499.25
$ws.Cells.Item(22, 12).Value = 5331.6665
$ws.Cells.Item(22, 13).Value = -149.25
$ws.Cells.Item(22, 14).Value = -6031.6665
$ws.Cells.Item(100, 8).Value = 70798
$ws.Cells.Item(100, 10).Value = 70798
$ws.Cells.Item(100, 12).Value = 70798
$ws.Cells.Item(100, 14).Value = -72962
$ws.Cells.Item(134, 8).Value = 2473.2593
$ws.Cells.Item(134, 9).Value = 1950.64
$ws.Cells.Item(134, 11).Value = 5851.92
$ws.Cells.Item(134, 13).Value = -3316.92
$ws.Cells.Item(141, 8).Value = 251443.5
$ws.Cells.Item(141, 10).Value = 251443.5
$ws.Cells.Item(141, 12).Value = 251443.5
$ws.Cells.Item(141, 14).Value = -261803.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 143.3077
$ws.Cells.Item(2, 10).Value = 177.36363
$ws.Cells.Item(2, 12).Value = 1064.18178
$ws.Cells.Item(2, 14).Value = -1290.18178
$ws.Cells.Item(32, 8).Value = 33983.332
$ws.Cells.Item(32, 10).Value = 33983.332
$ws.Cells.Item(32, 12).Value = 101949.996
$ws.Cells.Item(32, 14).Value = -102515.996
$ws.Cells.Item(40, 8).Value = 77.333336
$ws.Cells.Item(40, 10).Value = 120
$ws.Cells.Item(40, 12).Value = 480
$ws.Cells.Item(40, 14).Value = -618
$ws.Cells.Item(117, 8).Value = 4564.857
$ws.Cells.Item(117, 9).Value = 3000
$ws.Cells.Item(117, 10).Value = 4825.6665
$ws.Cells.Item(117, 11).Value = 9000
$ws.Cells.Item(117, 12).Value = 14476.9995
$ws.Cells.Item(117, 13).Value = -5558
$ws.Cells.Item(117, 14).Value = -21360.9995
$ws.Cells.Item(124, 8).Value = 18519986
$ws.Cells.Item(124, 9).Value = 746
$ws.Cells.Item(124, 10).Value = 25642770
$ws.Cells.Item(124, 11).Value = 2238
$ws.Cells.Item(124, 12).Value = 76928310
$ws.Cells.Item(124, 13).Value = 2672
$ws.Cells.Item(124, 14).Value = -76938130
$ws.Cells.Item(140, 8).Value = 2113.6
$ws.Cells.Item(140, 9).Value = 2113.6
$ws.Cells.Item(140, 11).Value = 6340.799999999999
$ws.Cells.Item(140, 13).Value = -1160.799999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(39, 8).Value = 90000
$ws.Cells.Item(39, 10).Value = 90000
$ws.Cells.Item(39, 12).Value = 90000
$ws.Cells.Item(39, 14).Value = -91064
$ws.Cells.Item(97, 8).Value = 1954.1538
$ws.Cells.Item(97, 9).Value = 1488.2222
$ws.Cells.Item(97, 11).Value = 1488.2222
$ws.Cells.Item(97, 13).Value = -992.2221999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 15531.576
$ws.Cells.Item(61, 9).Value = 17974.46
$ws.Cells.Item(61, 10).Value = 6458
$ws.Cells.Item(61, 11).Value = 17974.46
$ws.Cells.Item(61, 12).Value = 6458
$ws.Cells.Item(61, 13).Value = -17772.46
$ws.Cells.Item(61, 14).Value = -6862
$ws.Cells.Item(100, 8).Value = 12602.4
$ws.Cells.Item(100, 10).Value = 13288.286
$ws.Cells.Item(100, 12).Value = 13288.286
$ws.Cells.Item(100, 14).Value = -14370.286
$ws.Cells.Item(113, 8).Value = 15531.576
$ws.Cells.Item(113, 9).Value = 17974.46
$ws.Cells.Item(113, 10).Value = 6458
$ws.Cells.Item(113, 11).Value = 17974.46
$ws.Cells.Item(113, 12).Value = 6458
$ws.Cells.Item(113, 13).Value = -15804.46
$ws.Cells.Item(113, 14).Value = -10798
$ws.Cells.Item(136, 8).Value = 8754
$ws.Cells.Item(136, 9).Value = 6120
$ws.Cells.Item(136, 10).Value = 14022
$ws.Cells.Item(136, 11).Value = 18360
$ws.Cells.Item(136, 12).Value = 42066
$ws.Cells.Item(136, 13).Value = -15810
$ws.Cells.Item(136, 14).Value = -47166

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 4539.4165
$ws.Cells.Item(81, 9).Value = 2134.2
$ws.Cells.Item(81, 11).Value = 4268.4
$ws.Cells.Item(81, 13).Value = -3207.4
$ws.Cells.Item(84, 8).Value = 4539.4165
$ws.Cells.Item(84, 9).Value = 2134.2
$ws.Cells.Item(84, 11).Value = 21342
$ws.Cells.Item(84, 13).Value = -16038
$ws.Cells.Item(132, 8).Value = 5817.9165
$ws.Cells.Item(132, 9).Value = 3001.2222
$ws.Cells.Item(132, 11).Value = 9003.6666
$ws.Cells.Item(132, 13).Value = -6473.6666
$ws.Cells.Item(136, 8).Value = 4375.8
$ws.Cells.Item(136, 9).Value = 2679.1428
$ws.Cells.Item(136, 11).Value = 8037.428400000001
$ws.Cells.Item(136, 13).Value = -5487.428400000001

